# Add data for 2022-10-03
# - Renames the sheet / updates the "through <date>" labels from
#   September 23 to September 25.
# - Bumps the carjacking counts for the (partial) September 2022 column
#   ("B" .. "BD", shared string in B1) to reflect the newly-added data,
#   including a handful of neighborhoods that previously had no
#   carjackings recorded for that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab name and the workbook title both carry the "through" date.
$ws.Name = "Through 2022-09-25"

# B1 holds the shared-string column header for the current (partial)
# September 2022 month; update its text in place.
$ws.Range("B1").Value = "September 2022 (through September 25)"

# Updated counts: existing cells whose value increased, and cells that
# previously had no carjacking recorded (now getting their first count).
$updates = @{
    "K2"   = 10
    "AU2"  = 7
    "BD2"  = 2
    "K3"   = 15
    "T3"   = 6
    "AC3"  = 5
    "AU3"  = 6
    "BD3"  = 2
    "B4"   = 2
    "AU4"  = 2
    "K5"   = 3
    "B6"   = 6
    "K7"   = 5
    "T7"   = 3
    "AC7"  = 2
    "K8"   = 5
    "K10"  = 10
    "T10"  = 8
    "B11"  = 3
    "K11"  = 1
    "B12"  = 4
    "T12"  = 10
    "AL14" = 5
    "AU17" = 2
    "B19"  = 2
    "K20"  = 1
    "AL22" = 2
    "AU24" = 4
    "T25"  = 2
    "B26"  = 1
    "K26"  = 2
    "T33"  = 6
    "B34"  = 2
    "K38"  = 5
    "T38"  = 1
    "K39"  = 2
    "K45"  = 1
    "T46"  = 1
    "AU48" = 1
    "K50"  = 2
    "B54"  = 1
    "AC54" = 2
    "K58"  = 1
    "K63"  = 3
    "B64"  = 3
    "K77"  = 4
    "AC89" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
